# outputs-HGR-r202/test-g__CAG-631_split_pruned.xlsx
#
# - B2 / B3 get new numeric values (previously both were literal 1s, a
#   placeholder/"pruned" copy; the real prediction scores are restored).
# - The header row (A1:C1) plus the two "Row" label cells (A2/A3) get
#   re-stamped with a fresh (but format-identical: text / "@") cell style,
#   matching how the upstream ful-path.csv -> xlsx exporter re-wrote the
#   style table on this pass.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-stamp the label column / header styles with a new style-table entry
# (same text number format as before) by toggling a format property off
# and back on - this mints a fresh style index instead of reusing the old
# one, mirroring the re-export the diff captures.
foreach ($addr in @("A1", "B1", "C1", "A2", "A3")) {
    $cell = $ws.Range($addr)
    $cell.Locked = $false
    $cell.Locked = $true
}

# Restore the real prediction values (previously placeholder 1s).
$ws.Range("B2").Value = 0.087143138009892596
$ws.Range("B3").Value = -0.04862773867910164
